$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header ---
$ws.Range("C1").Value = "Annotated"

# --- Team rows 2..31: (row, GameCount, Annotated-or-null) ---
$data = @(
    @(2,  4, 1),
    @(3,  3, 1),
    @(4,  4, 1),
    @(5,  4, $null),
    @(6,  4, 3),
    @(7,  5, 1),
    @(8,  2, $null),
    @(9,  7, $null),
    @(10, 8, 3),
    @(11, 5, $null),
    @(12, 3, $null),
    @(13, 4, 1),
    @(14, 5, 2),
    @(15, 6, 4),
    @(16, 2, $null),
    @(17, 4, 2),
    @(18, 4, 1),
    @(19, 4, $null),
    @(20, 5, $null),
    @(21, 6, 2),
    @(22, 4, $null),
    @(23, 3, $null),
    @(24, 3, $null),
    @(25, 4, $null),
    @(26, 5, 1),
    @(27, 4, 2),
    @(28, 1, 2),
    @(29, 4, $null),
    @(30, 4, 2),
    @(31, 3, 1)
)

foreach ($row in $data) {
    $r = $row[0]
    $gameCount = $row[1]
    $annotated = $row[2]

    $ws.Cells.Item($r, 2).Value = $gameCount
    if ($annotated -ne $null) {
        $ws.Cells.Item($r, 3).Value = $annotated
    }
}

# --- Summary rows 33..36 get column C formulas too ---
$ws.Range("C33").Formula = "=SUM(C2:C31)/2"

$ws.Range("C34").Formula = "=AVERAGE(C2:C31)"
$ws.Range("C34").NumberFormat = "0.0"

$ws.Range("C35").Formula = "=MAX(C2:C31)"

$ws.Range("C36").Formula = "=MIN(C2:C31)"

# --- Cosmetic: match the author's final cell selection ---
[void]$ws.Range("E19").Select()

